$wb = $excel.ActiveWorkbook

# --- weibull ---
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -3.01920941114309
$ws.Range("C2").Value = 0.1957543823296
$ws.Range("B3").Value = 0.0688710722139505
$ws.Range("C3").Value = 0.106918285891009

# --- lognormal ---
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 2.46094374220508
$ws.Range("C2").Value = 0.264036571292589
$ws.Range("B3").Value = -1.01546611467525
$ws.Range("C3").Value = 0.109562714690895

# --- llogis ---
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -2.40113898532518
$ws.Range("C2").Value = 0.122189329723341
$ws.Range("B3").Value = 0.565964498875003
$ws.Range("C3").Value = 0.112695489584321

# --- gompertz ---
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.64863060347219
$ws.Range("C2").Value = 0.140374272993617
$ws.Range("B3").Value = -0.0168393407736236
$ws.Range("C3").Value = 0.0119133252065979

# --- weibull cov ---
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.0383197782012431
$ws.Range("B2").Value = -0.0153403431417579
$ws.Range("A3").Value = -0.0153403431417579
$ws.Range("B3").Value = 0.0114315198578715

# --- lognormal cov ---
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.0697153109799464
$ws.Range("B2").Value = -0.0254842257230339
$ws.Range("A3").Value = -0.0254842257230339
$ws.Range("B3").Value = 0.0120039884504384

# --- llogis cov ---
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.0149302322982393
$ws.Range("B2").Value = 0.00300575235300295
$ws.Range("A3").Value = 0.00300575235300295
$ws.Range("B3").Value = 0.0127002733726498

# --- gompertz cov ---
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0197049365184865
$ws.Range("B2").Value = -0.000937601013814208
$ws.Range("A3").Value = -0.000937601013814208
$ws.Range("B3").Value = 0.000141927317478162
